$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" '244.49'
Set-TextValue "D3" '21.76'
Set-TextValue "D4" '5.392'
Set-TextValue "D5" '0.06011'
Set-TextValue "D6" '3.390'
Set-TextValue "D7" '0.8153'
Set-TextValue "D8" '0.9473'
Set-TextValue "D9" '0.1436'
Set-TextValue "D10" '0.07425'
Set-TextValue "D11" '0.03329'
Set-TextValue "D12" '0.03057'
Set-TextValue "D13" '0.09408'
Set-TextValue "D14" '4.001'
Set-TextValue "D16" '0.04817'
Set-TextValue "D17" '0.0005911'
Set-TextValue "D18" '0.005511'
Set-TextValue "D19" '0.004162'
Set-TextValue "D20" '0.0009865'
Set-TextValue "D22" '6.423'
Set-TextValue "D26" '0.00007003'
Set-TextValue "D40" '0.04019'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D41" '0.006412'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1073'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002901'
$ws.Range("E43").Value = '42CEJICEJI'
Set-TextValue "D44" '0.006580'
Set-TextValue "D45" '0.00005250'
Set-TextValue "D47" '0.8602'
Set-TextValue "D48" '0.003261'
